$wb = $excel.ActiveWorkbook

# --- Summary sheet: move selection only (no data change) ---
$wsSummary = $wb.Worksheets.Item("Summary")
[void]$wsSummary.Range("D4").Select()

# --- Repayment schedule sheet: remove the duplicated "Over Due" (column O) values ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
[void]$wsSchedule.Range("O3:O14").Clear()
[void]$wsSchedule.Range("P2").Clear()
[void]$wsSchedule.Range("F11").Select()

# --- Transactions sheet: update transaction IDs and move selection ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 1943
$wsTransactions.Range("A3").Value = 1942
$wsTransactions.Range("A4").Value = 1941
$wsTransactions.Range("A5").Value = 1940
[void]$wsTransactions.Range("D5").Select()
